# "Generate Report for Handoff"
# Updates status from "In Translation" to "Ready for handoff" on all sheets,
# and refreshes the handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status columns: "In Translation" -> "Ready for handoff"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Latest Handoff Date / Datetime refreshed to reflect new handoff
$wsOverview.Range("D2").Value = "2016-30-12 04:30:55"
$wsZhCn.Range("E2").Value = "2016-03-12 04:30:53"
$wsDeDe.Range("E2").Value = "2016-03-12 04:30:55"
